# Update the "Förändrad" (Changed) date column C for all data rows
# from serial date 45189 (2023-09-20) to 45190 (2023-09-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 540  # data occupies rows 2 through 540

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
